$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# res_bus/vm_pu.xlsx -- case with 380 kV done: slack bus voltage setpoint
# (bus 0 / column B) lowered from 1.05 to 1.02 p.u., and the resulting
# power-flow solution for every other bus voltage magnitude is updated
# for all 24 time steps (rows 2-25).

# row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039432910895759
$ws.Range("D2").Value = 1.045305475662108
$ws.Range("E2").Value = 1.043087205207078
$ws.Range("F2").Value = 1.05295446749078
$ws.Range("I2").Value = 1.039910959385147
$ws.Range("J2").Value = 1.044524955302397
$ws.Range("K2").Value = 1.048074088137139
$ws.Range("L2").Value = 1.045862060540546
$ws.Range("M2").Value = 1.055701772141627
$ws.Range("N2").Value = 1.046008300547614

# row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040349853320488
$ws.Range("D3").Value = 1.046019592391631
$ws.Range("E3").Value = 1.043952670828574
$ws.Range("F3").Value = 1.05385442162665
$ws.Range("I3").Value = 1.040130142079999
$ws.Range("J3").Value = 1.045087264202358
$ws.Range("K3").Value = 1.04859987197324
$ws.Range("L3").Value = 1.046538347188166
$ws.Range("M3").Value = 1.056414451656517
$ws.Range("N3").Value = 1.04657140799071

# row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04094373301865
$ws.Range("D4").Value = 1.046482163576676
$ws.Range("E4").Value = 1.044513574766649
$ws.Range("F4").Value = 1.054437705180361
$ws.Range("I4").Value = 1.040271053598669
$ws.Range("J4").Value = 1.045451035189886
$ws.Range("K4").Value = 1.048939905810231
$ws.Range("L4").Value = 1.046976198462143
$ws.Range("M4").Value = 1.056875919927471
$ws.Range("N4").Value = 1.046935695574766

# row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041193532113462
$ws.Range("D5").Value = 1.046676744305646
$ws.Range("E5").Value = 1.044749590310581
$ws.Range("F5").Value = 1.054683144185353
$ws.Range("I5").Value = 1.040330073317881
$ws.Range("J5").Value = 1.045603944188256
$ws.Range("K5").Value = 1.049082811028869
$ws.Range("L5").Value = 1.047160329420044
$ws.Range("M5").Value = 1.057069995746915
$ws.Range("N5").Value = 1.047088821721452

# row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041235482190803
$ws.Range("D6").Value = 1.046709422013871
$ws.Range("E6").Value = 1.044789230719667
$ws.Range("F6").Value = 1.054724367716756
$ws.Range("I6").Value = 1.040339970100786
$ws.Range("J6").Value = 1.045629617059289
$ws.Range("K6").Value = 1.049106802773183
$ws.Range("L6").Value = 1.047191249200052
$ws.Range("M6").Value = 1.057102586264639
$ws.Range("N6").Value = 1.047114531050907

# row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040947070328915
$ws.Range("D7").Value = 1.046484763120733
$ws.Range("E7").Value = 1.044516727588606
$ws.Range("F7").Value = 1.054440983859935
$ws.Range("I7").Value = 1.040271843086195
$ws.Range("J7").Value = 1.045453078448593
$ws.Range("K7").Value = 1.048941815494958
$ws.Range("L7").Value = 1.046978658601543
$ws.Range("M7").Value = 1.056878512886762
$ws.Range("N7").Value = 1.046937741735134

# row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03974267990948
$ws.Range("D8").Value = 1.045546712423987
$ws.Range("E8").Value = 1.043379508226576
$ws.Range("F8").Value = 1.053258412929061
$ws.Range("I8").Value = 1.039985222169234
$ws.Range("J8").Value = 1.044715006252269
$ws.Range("K8").Value = 1.048251816580774
$ws.Range("L8").Value = 1.046090562566523
$ws.Range("M8").Value = 1.055942559092838
$ws.Range("N8").Value = 1.046198621391634

# row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.037624705934727
$ws.Range("D9").Value = 1.043897564879569
$ws.Range("E9").Value = 1.041382463331824
$ws.Range("F9").Value = 1.051181940199761
$ws.Range("I9").Value = 1.039473182783698
$ws.Range("J9").Value = 1.043413859288899
$ws.Range("K9").Value = 1.047034593199967
$ws.Range("L9").Value = 1.044527582218009
$ws.Range("M9").Value = 1.054295775113598
$ws.Range("N9").Value = 1.044895626650332

# row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036215695154769
$ws.Range("D10").Value = 1.042800787987931
$ws.Range("E10").Value = 1.040055815003823
$ws.Range("F10").Value = 1.049802671863046
$ws.Range("I10").Value = 1.039127167384224
$ws.Range("J10").Value = 1.042546107887317
$ws.Range("K10").Value = 1.046222262832804
$ws.Range("L10").Value = 1.043486984942126
$ws.Range("M10").Value = 1.053199669284627
$ws.Range("N10").Value = 1.044026642942227

# row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035606297329724
$ws.Range("D11").Value = 1.042326519648274
$ws.Range("E11").Value = 1.039482498491985
$ws.Range("F11").Value = 1.049206649908066
$ws.Range("I11").Value = 1.038976241889744
$ws.Range("J11").Value = 1.042170299129316
$ws.Range("K11").Value = 1.045870327110386
$ws.Range("L11").Value = 1.043036739676363
$ws.Range("M11").Value = 1.052725474914016
$ws.Range("N11").Value = 1.043650300492684

# row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035380048117596
$ws.Range("D12").Value = 1.042150453397509
$ws.Range("E12").Value = 1.039269714608582
$ws.Range("F12").Value = 1.048985444208391
$ws.Range("I12").Value = 1.038920016961678
$ws.Range("J12").Value = 1.042030698040169
$ws.Range("K12").Value = 1.045739574834906
$ws.Range("L12").Value = 1.042869550860836
$ws.Range("M12").Value = 1.052549403423109
$ws.Range("N12").Value = 1.043510501153978

# row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035428574423778
$ws.Range("D13").Value = 1.042188215741737
$ws.Range("E13").Value = 1.039315349694061
$ws.Range("F13").Value = 1.049032885258032
$ws.Range("I13").Value = 1.038932084833564
$ws.Range("J13").Value = 1.042060643337586
$ws.Range("K13").Value = 1.045767622877239
$ws.Range("L13").Value = 1.042905411042131
$ws.Range("M13").Value = 1.052587168374569
$ws.Range("N13").Value = 1.043540488977152

# row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035587593283954
$ws.Range("D14").Value = 1.042311963946733
$ws.Range("E14").Value = 1.039464906204808
$ws.Range("F14").Value = 1.04918836122579
$ws.Range("I14").Value = 1.038971597671708
$ws.Range("J14").Value = 1.042158759837014
$ws.Range("K14").Value = 1.045859519644972
$ws.Range("L14").Value = 1.043022918722337
$ws.Range("M14").Value = 1.052710919447704
$ws.Range("N14").Value = 1.043638744813263

# row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035685584433254
$ws.Range("D15").Value = 1.042388222347138
$ws.Range("E15").Value = 1.039557075663026
$ws.Range("F15").Value = 1.049284179448623
$ws.Range("I15").Value = 1.038995921057177
$ws.Range("J15").Value = 1.042219211491113
$ws.Range("K15").Value = 1.045916136640781
$ws.Range("L15").Value = 1.043095326051902
$ws.Range("M15").Value = 1.052787175274998
$ws.Range("N15").Value = 1.043699282315645

# row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036256154002982
$ws.Range("D16").Value = 1.042832277295427
$ws.Range("E16").Value = 1.040093888087106
$ws.Range("F16").Value = 1.049842253510154
$ws.Range("I16").Value = 1.039137160715291
$ws.Range("J16").Value = 1.042571047760305
$ws.Range("K16").Value = 1.046245615707699
$ws.Range("L16").Value = 1.043516873508789
$ws.Range("M16").Value = 1.053231149104597
$ws.Range("N16").Value = 1.044051618232696

# row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036614248945385
$ws.Range("D17").Value = 1.043110994505087
$ws.Range("E17").Value = 1.040430920175554
$ws.Range("F17").Value = 1.050192643761945
$ws.Range("I17").Value = 1.039225462907619
$ws.Range("J17").Value = 1.04279172823486
$ws.Range("K17").Value = 1.046452238752638
$ws.Range("L17").Value = 1.043781391031142
$ws.Range("M17").Value = 1.053509757316436
$ws.Range("N17").Value = 1.044272612098836

# row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036823188398317
$ws.Range("D18").Value = 1.043273627507065
$ws.Range("E18").Value = 1.040627614313304
$ws.Range("F18").Value = 1.050397137143368
$ws.Range("I18").Value = 1.03927686205098
$ws.Range("J18").Value = 1.042920440919469
$ws.Range("K18").Value = 1.046572739948579
$ws.Range("L18").Value = 1.043935712371921
$ws.Range("M18").Value = 1.053672305854028
$ws.Range("N18").Value = 1.044401507570221

# row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036894442975979
$ws.Range("D19").Value = 1.043329091623932
$ws.Range("E19").Value = 1.040694700378358
$ws.Range("F19").Value = 1.050466883844854
$ws.Range("I19").Value = 1.03929436982869
$ws.Range("J19").Value = 1.042964327484319
$ws.Range("K19").Value = 1.046613824566931
$ws.Range("L19").Value = 1.043988337478521
$ws.Range("M19").Value = 1.053727737633083
$ws.Range("N19").Value = 1.044445456459028

# row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036575821629625
$ws.Range("D20").Value = 1.043081084365797
$ws.Range("E20").Value = 1.040394748534675
$ws.Range("F20").Value = 1.050155038124454
$ws.Range("I20").Value = 1.039215999881322
$ws.Range("J20").Value = 1.04276805197356
$ws.Range("K20").Value = 1.046430071964369
$ws.Range("L20").Value = 1.043753007408628
$ws.Range("M20").Value = 1.05347986103335
$ws.Range("N20").Value = 1.04424890221453

# row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035540763195073
$ws.Range("D21").Value = 1.042275520464023
$ws.Range("E21").Value = 1.039420860796723
$ws.Range("F21").Value = 1.049142572353417
$ws.Range("I21").Value = 1.0389599666616
$ws.Range("J21").Value = 1.042129867213101
$ws.Range("K21").Value = 1.045832459092841
$ws.Range("L21").Value = 1.04298831418522
$ws.Range("M21").Value = 1.052674476024887
$ws.Range("N21").Value = 1.04360981115851

# row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034890607877822
$ws.Range("D22").Value = 1.041769598933639
$ws.Range("E22").Value = 1.038809531301152
$ws.Range("F22").Value = 1.048507056759879
$ws.Range("I22").Value = 1.038798037068248
$ws.Range("J22").Value = 1.041728563661911
$ws.Range("K22").Value = 1.045456556551681
$ws.Range("L22").Value = 1.042507824531738
$ws.Range("M22").Value = 1.052168477605636
$ws.Range("N22").Value = 1.043207937710247

# row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035235207700756
$ws.Range("D23").Value = 1.042037743130911
$ws.Range("E23").Value = 1.039133514072534
$ws.Range("F23").Value = 1.048843854432909
$ws.Range("I23").Value = 1.038883968997018
$ws.Range("J23").Value = 1.041941306799709
$ws.Range("K23").Value = 1.04565584437437
$ws.Range("L23").Value = 1.042762512032154
$ws.Range("M23").Value = 1.052436680570117
$ws.Range("N23").Value = 1.043420982967703

# row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036593185067262
$ws.Range("D24").Value = 1.043094599278972
$ws.Range("E24").Value = 1.040411092605501
$ws.Range("F24").Value = 1.050172030132985
$ws.Range("I24").Value = 1.039220276143251
$ws.Range("J24").Value = 1.042778750277107
$ws.Range("K24").Value = 1.046440088239199
$ws.Range("L24").Value = 1.043765832644835
$ws.Range("M24").Value = 1.053493369750013
$ws.Range("N24").Value = 1.044259615710895

# row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038171734318072
$ws.Range("D25").Value = 1.044323447568319
$ws.Range("E25").Value = 1.041897922759267
$ws.Range("F25").Value = 1.051857541386154
$ws.Range("I25").Value = 1.039606380423832
$ws.Range("J25").Value = 1.043750297989019
$ws.Range("K25").Value = 1.047349428823373
$ws.Range("L25").Value = 1.044931410519958
$ws.Range("M25").Value = 1.054721205528858
$ws.Range("N25").Value = 1.045232543131996

Write-Output "updated vm_pu values for rows 2-25"